$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab19")

# Fix mangled accented characters in the Regional Economic Communities footnote (cell A103)
$ws.Range("A103").Value = 'Regional Economic Communities:CEN-SAD = "Community of Sahel-Saharan States";COMESA = "Common Market for Eastern and Southern Africa";EAC = "East African Community";ECCAS = "Economic Community of Central African States";ECOWAS = "Economic Community of West African States";IGAD = "Intergovernmental Authority on Development";SADC = "Southern African Development Community";UMA = "Arab Maghreb Union";PALOP = "Países Africanos de Língua Oficial Portuguesa";ASEAN = "Association of Southeast Asian Nations";MERCOSUR = "Mercado Común del Sur".EU27 = "European Union (27 members)".OECD = "Organisation for Economic Co-operation and Development".'

# Minor recalculated value updates (row 67)
$ws.Range("H67").Value = 5.29175561613473
$ws.Range("L67").Value = 18.6011441753639

# Minor recalculated value update (row 70)
$ws.Range("N70").Value = 3.00757906215029

# Recalculated values (row 97)
$ws.Range("C97").Value = 9.23299085052547
$ws.Range("D97").Value = 7.11636757045701
$ws.Range("E97").Value = 11.3709537855184
$ws.Range("F97").Value = 15.8363945094399
$ws.Range("G97").Value = 24.4636653418919
$ws.Range("H97").Value = 4.46694613562718
$ws.Range("I97").Value = 11.8078131045872
$ws.Range("J97").Value = 9.77502805063377
$ws.Range("K97").Value = 13.507604281606
$ws.Range("L97").Value = 25.9832987793666
$ws.Range("M97").Value = 18.0899481753386
$ws.Range("N97").Value = 3.69109433978058
$ws.Range("O97").Value = 21583.432691000002
$ws.Range("P97").Value = 96378.972787999999
$ws.Range("Q97").Value = 18.2968739941831
$ws.Range("R97").Value = 18.5436913281488

# Recalculated values (row 98)
$ws.Range("C98").Value = 6.58839387364423
$ws.Range("E98").Value = 42.9993599616683
$ws.Range("F98").Value = 14.5464957020326
$ws.Range("G98").Value = 23.9815661593687
$ws.Range("H98").Value = 18.1395408192066
$ws.Range("I98").Value = 11.9716161625933
$ws.Range("K98").Value = 74.1020277176156
$ws.Range("L98").Value = 28.0554149027077
$ws.Range("M98").Value = 6.47958422010887
$ws.Range("N98").Value = 7.45461227489813
$ws.Range("O98").Value = 153212.318306
$ws.Range("P98").Value = 54680.485189000101
$ws.Range("Q98").Value = 73.697749864480002
$ws.Range("R98").Value = 19.2558164737962
